$wb = $excel.ActiveWorkbook

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")

# "cb985dd8" row (row 7) got a new (later) handoff timestamp for zh-cn.
# The same timestamp string is shared by rows 10, 11, 14, 15, 16 ("Ready for handoff"
# rows that still carry the earlier handoff run), so they are refreshed together
# to keep pointing at one shared value.
$wsZh.Range("D7").Value = "2016-03-10 18:24:36"
$wsZh.Range("D10").Value = "2016-03-10 18:24:36"
$wsZh.Range("D11").Value = "2016-03-10 18:24:36"
$wsZh.Range("D14").Value = "2016-03-10 18:24:36"
$wsZh.Range("D15").Value = "2016-03-10 18:24:36"
$wsZh.Range("D16").Value = "2016-03-10 18:24:36"

# "77a29242" file (rows 12 & 13) now reuses the same handoff timestamp instead of
# keeping its own separate (now stale/duplicate) one.
$wsZh.Range("D12").Value = "2016-03-10 18:24:36"
$wsZh.Range("D13").Value = "2016-03-10 18:24:36"

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")

# Same regeneration applied to the de-de report.
$wsDe.Range("D7").Value = "2016-03-10 18:24:42"
$wsDe.Range("D10").Value = "2016-03-10 18:24:42"
$wsDe.Range("D11").Value = "2016-03-10 18:24:42"
$wsDe.Range("D14").Value = "2016-03-10 18:24:42"
$wsDe.Range("D15").Value = "2016-03-10 18:24:42"
$wsDe.Range("D16").Value = "2016-03-10 18:24:42"

$wsDe.Range("D12").Value = "2016-03-10 18:24:42"
$wsDe.Range("D13").Value = "2016-03-10 18:24:42"
